$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to start with a plain row-number column (A) followed by
# the real headings in B:F. Turn this into a proper table starting at A1
# by dropping the row-number column - this shifts B:F left into A:E.
$ws.Columns.Item(1).Delete()

# Dynamically add two new columns to the table and fill in their heading.
$ws.Range("F1").Value = "Test1"
$ws.Range("G1").Value = "Test2"

# Add the data for the new columns below the headings.
$ws.Range("F2").Value = 12345
$ws.Range("G2").Value = 4567
$ws.Range("F3").Value = 8787
$ws.Range("G3").Value = 999

# The new columns are fresh/default formatted, not part of the original
# centered table styling.
$ws.Range("F1:G3").HorizontalAlignment = 1

# Keep the original table (now A1:E3) consistently center-aligned.
$ws.Range("A1:E3").HorizontalAlignment = -4108

# Resize the table's columns to fit its new, smaller, dynamic content.
$ws.Columns.Item(1).ColumnWidth = 11.58
$ws.Columns.Item(2).ColumnWidth = 4.76
$ws.Columns.Item(3).ColumnWidth = 7.54
$ws.Columns.Item(4).ColumnWidth = 6.71
$ws.Columns.Item(5).ColumnWidth = 6.71

# Leave the selection where the editor ended up.
[void]$ws.Range("I7").Select()
